# Apply the "purple (7030A0) Q&A highlight" edit:
#   - Color the "Part 1" / "Part 2" headings in the "Your Answers" section purple.
#   - Bold + color the "Looking at all of the yammer datasets..." list item, add a
#     blank paragraph right after it.
#   - Color the "Looking at the yammer 'publisher update'..." list item, add a new
#     indented blank paragraph right after it.

$d = $word.ActiveDocument
$purple = 10498160   # RGB(0x70,0x30,0xA0) -> 0x70 + 0x30*256 + 0xA0*65536 = 7030A0 in Word's BGR long form

$wordXmlOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$wordXmlClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---- locate the target paragraphs --------------------------------------
# The doc contains two copies of this Q&A block: an earlier "instructions"
# copy and the real "Your Answers" copy near the end. We want the LAST
# occurrence of each target paragraph (the "Your Answers" one) - the loop
# below keeps overwriting the index on every match, so it naturally ends
# up on the final (bottom-most) occurrence in the document.

$part1Index = -1
$part2Index = -1
$q1Index = -1
$q2Index = -1

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text

    if (($t -eq "Part 1`r") -or ($t -eq "Part 1")) {
        if ($p.Style.NameLocal -eq "Normal") {
            $part1Index = $idx
        }
    }
    if (($t -eq "Part 2`r") -or ($t -eq "Part 2")) {
        if ($p.Style.NameLocal -eq "Heading 2") {
            $part2Index = $idx
        }
    }
    if ($t -like "Looking at all of the yammer datasets*") {
        $q1Index = $idx
    }
    if ($t -like "Looking at the yammer*publisher update*") {
        $q2Index = $idx
    }
}

Write-Host "Part1:" $part1Index "Part2:" $part2Index "Q1:" $q1Index "Q2:" $q2Index

if (($part1Index -eq -1) -or ($part2Index -eq -1) -or ($q1Index -eq -1) -or ($q2Index -eq -1)) {
    throw "Could not locate all target paragraphs for the edit."
}

# ---- 1) "Part 1" heading (bold/sz28 run in the Your Answers section) ---
$p1 = $d.Paragraphs.Item($part1Index)
$p1.Range.Font.Color = $purple

# ---- 2) "Part 2" heading (Heading 2 style) ------------------------------
$p2 = $d.Paragraphs.Item($part2Index)
$p2.Range.Font.Color = $purple

# ---- 3) First list question: bold + purple, then a blank <w:p/> after --
# Use InsertXML (full paragraph replace-in-place) rather than Font.Bold /
# Font.BoldBi, because the latter only ever stamps <w:b/> (no <w:bCs/>) on
# the paragraph-mark's own run properties (<w:pPr><w:rPr>) in this host,
# even though it correctly stamps <w:b/><w:bCs/> on the text run itself.
# InsertXML lets us state the exact target markup for both places.
$q1 = $d.Paragraphs.Item($q1Index)
$q1Xml = $wordXmlOpen + '<w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:b/><w:bCs/><w:color w:val="7030A0"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="7030A0"/></w:rPr><w:t>Looking at all of the yammer datasets, what conclusions can you draw about the types of users who have the highest engagements with Yammer?</w:t></w:r></w:p></w:body>' + $wordXmlClose
$q1.Range.InsertXML($q1Xml)

# Insert a bare empty paragraph right after it.
$q1 = $d.Paragraphs.Item($q1Index)
$q1.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item($q1Index + 1)
$blankXml = $wordXmlOpen + '<w:body><w:p/></w:body>' + $wordXmlClose
$blank1.Range.InsertXML($blankXml)

# ---- 4) Second list question: purple only, then an indented blank after
$q2Index = $q2Index + 1   # shifted by the blank paragraph inserted above
$q2 = $d.Paragraphs.Item($q2Index)
$q2.Range.Font.Color = $purple

$q2 = $d.Paragraphs.Item($q2Index)
$q2.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item($q2Index + 1)
$indentXml = $wordXmlOpen + '<w:body><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p></w:body>' + $wordXmlClose
$blank2.Range.InsertXML($indentXml)

Write-Host "Applied edits."
